# Generate Report for Handback
# Update the localization status workbook: mark the 81d9b0b3... file as
# handed back (in sync with en-US) across all sheets, refresh the
# "Latest Handback DateTime" / clear the stale error detail for the
# zh-cn and de-de language sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the 81d9b0b3... file ---
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row 3 is the 81d9b0b3... file ---
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-08-31 02:55:33"
$wsZhCn.Range("P3").Value = ""

# --- de-de sheet: row 3 is the 81d9b0b3... file ---
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-08-31 02:55:41"
$wsDeDe.Range("P3").Value = ""

# Error Detail column is no longer needed now that handback succeeded,
# narrow it back down from its error-message width.
$wsZhCn.Columns.Item(16).ColumnWidth = 13.7470528738839
$wsDeDe.Columns.Item(16).ColumnWidth = 13.7470528738839
